$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2: constant used by the IF formulas below (Aile yardımı = 100)
$ws.Range("E2").Value = 100

# E4: standalone IF formula
$ws.Range("E4").Formula = "=IF(D4=""HAYIR"",C4+`$E`$2,C4)"

# E5:E11: same IF formula filled down as one shared formula group
$ws.Range("E5:E11").Formula = "=IF(D5=""HAYIR"",C5+`$E`$2,C5)"

# Student info block (Numara / Ad Soyad / Bölüm)
$ws.Range("J7").Value = 20215070019
$ws.Range("J8").Value = "KÜBRA ÇABUK"
$ws.Range("J9").Value = "YBS"

# Update the active selection to match the saved view
$null = $ws.Range("H3").Select()
